# Weekly fruit/vegetable price update: insert a new daily record as the
# first row of the "Haba" price-history block and push the rest of the
# existing rows down by one (349 -> 350, 350 -> 351, ... 407 -> 408).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 349..407 down to 350..408, leaving a blank row 349 behind.
$ws.Rows("349:349").Insert()

# Populate the newly inserted row 349 with the latest price record.
$ws.Range("A349").Value = 9
$ws.Range("B349").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C349").Value = "Metropolitana"
$ws.Range("D349").Value = 45218
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112026
$ws.Range("G349").Value = "Haba"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 70
$ws.Range("K349").Value = 10000
$ws.Range("L349").Value = 12000
$ws.Range("M349").Value = 11000
$ws.Range("N349").Value = "$/saco 25 kilos"
$ws.Range("O349").Value = "Provincia de Melipilla"
$ws.Range("P349").Value = 440
$ws.Range("Q349").Value = 25
$ws.Range("R349").Value = "Hortaliza"
